$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SMARTseq2")

# Clear empty wells
$ws.Range("B9").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("M16").Value = ""

# Replace G1-G4 (row15 B:E) with water variants
$ws.Range("K11").Value = "Water"
$ws.Range("C15").Value = "water"
$ws.Range("B15").Value = "  water"
$ws.Range("E15").Value = "wat   er"
$ws.Range("D15").Value = "wAter   "

# Selection
$ws.Range("D15").Select()
